$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update probability values in column A
$ws.Range("A6").Value = 0.11
$ws.Range("A7").Value = 0.08
$ws.Range("A8").Value = 0.04
$ws.Range("A9").Value = 0.04
$ws.Range("A11").Value = 0.02

# Convert the numeric symbol ids in column B to letter labels
$ws.Range("B2").Value = "A"
$ws.Range("B3").Value = "B"
$ws.Range("B4").Value = "C"
$ws.Range("B5").Value = "D"
$ws.Range("B6").Value = "E"
$ws.Range("B7").Value = "F"
$ws.Range("B8").Value = "G"
$ws.Range("B9").Value = "H"
$ws.Range("B10").Value = "I"
$ws.Range("B11").Value = "J"
$ws.Range("B12").Value = "K"

# Center-align column B (header + data) and size the column
$ws.Range("B1:B12").HorizontalAlignment = -4108
$ws.Columns("B").ColumnWidth = 10

# Move the active selection
$ws.Range("C6").Select()
